# Generate Report for Handoff
# The "Status" for the report's only tracked file moves from "In Translation"
# to "Ready for handoff", and the handoff timestamps advance a few minutes.
# Widening the "Status" columns (Overview!E:F, zh-cn!C, de-de!C) so the new,
# longer status text fits without truncation.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---
$ws_overview.Range("E2").Value = "Ready for handoff"
$ws_overview.Range("F2").Value = "Ready for handoff"
$ws_zhcn.Range("C2").Value = "Ready for handoff"
$ws_dede.Range("C2").Value = "Ready for handoff"

# --- Latest HO Xliff Generate Date / de-de Latest Handoff Datetime ---
$ws_overview.Range("G2").Value = "2016-08-20 10:45:40"
$ws_dede.Range("H2").Value = "2016-08-20 10:45:40"

# --- zh-cn Latest Handoff Datetime ---
$ws_zhcn.Range("H2").Value = "2016-08-20 10:45:36"

# --- Widen the Status columns to fit the new text ---
$newStatusWidth = (103 / 6) - (5 / 6)
$ws_overview.Columns("E").ColumnWidth = $newStatusWidth
$ws_overview.Columns("F").ColumnWidth = $newStatusWidth
$ws_zhcn.Columns("C").ColumnWidth = $newStatusWidth
$ws_dede.Columns("C").ColumnWidth = $newStatusWidth
